# Helper: pull a paragraph's own OOXML (<w:p ...> ... </w:p>) out of the
# flat-OPC payload that Range.WordOpenXML hands back, stripping the
# synthetic w14:paraId/w14:textId the packager stamps on for the roundtrip
# (the source document never used them).
function Get-ParaXml($para) {
    $full = $para.Range.WordOpenXML
    $bodyMatch = [regex]::Match($full, '<w:body>(.*)</w:body>', [System.Text.RegularExpressions.RegexOptions]::Singleline)
    $body = $bodyMatch.Groups[1].Value
    $pMatch = [regex]::Match($body, '(<w:p\b.*?</w:p>)', [System.Text.RegularExpressions.RegexOptions]::Singleline)
    $frag = $pMatch.Groups[1].Value
    $frag = $frag -replace ' w14:paraId="[0-9A-Fa-f]*"', ''
    $frag = $frag -replace ' w14:textId="[0-9A-Fa-f]*"', ''
    return $frag
}

function Remove-Highlight($xmlFrag) {
    return $xmlFrag -replace '<w:highlight[^/]*/>', ''
}

$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1. The stray "_GoBack" bookmark sitting right after the title line is
#    not wanted there any more -- drop it (it gets reinstated further
#    down, at the spot Word actually left the cursor).
# ---------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

# ---------------------------------------------------------------------
# 2. Locate the "Грантополучатель" signature block by content so the
#    script does not depend on brittle absolute paragraph indices.
# ---------------------------------------------------------------------
$idxLabel = -1
$i = 0
foreach ($p in $d.Paragraphs) {
    $i = $i + 1
    if ($p.Range.Text -like "*Грантополучатель*") {
        $idxLabel = $i
    }
}

function Get-ParaByIndex($index) {
    $i = 0
    foreach ($p in $d.Paragraphs) {
        $i = $i + 1
        if ($i -eq $index) { return $p }
    }
    return $null
}

$pLabel = Get-ParaByIndex $idxLabel
$pBlank = Get-ParaByIndex ($idxLabel + 1)
$pLine  = Get-ParaByIndex ($idxLabel + 2)
$pRoles = Get-ParaByIndex ($idxLabel + 3)

# ---------------------------------------------------------------------
# 3. "Грантополучатель" paragraph: strip the yellow highlight from both
#    the paragraph mark and the run.
# ---------------------------------------------------------------------
$frag = Get-ParaXml $pLabel
$frag = Remove-Highlight $frag
$pLabel.Range.InsertXML($frag)

# ---------------------------------------------------------------------
# 4. The blank paragraph right below it: same deal, mark-only highlight.
# ---------------------------------------------------------------------
$pBlank = Get-ParaByIndex ($idxLabel + 1)
$frag = Get-ParaXml $pBlank
$frag = Remove-Highlight $frag
$pBlank.Range.InsertXML($frag)

# ---------------------------------------------------------------------
# 5. The underscores/"/М.П/" line: drop the highlight and split the run
#    in two at the point Word's cursor ended up, re-planting the
#    "_GoBack" bookmark between the two halves.
# ---------------------------------------------------------------------
$pLine = Get-ParaByIndex ($idxLabel + 2)
$frag = Get-ParaXml $pLine
$frag = Remove-Highlight $frag

$runMatch = [regex]::Match($frag, '(<w:r\b[^>]*>)(<w:rPr>.*?</w:rPr>)(<w:t[^>]*>)(.*?)(</w:t></w:r>)', [System.Text.RegularExpressions.RegexOptions]::Singleline)
$runOpen   = $runMatch.Groups[1].Value
$runProps  = $runMatch.Groups[2].Value
$tOpen     = $runMatch.Groups[3].Value
$text      = $runMatch.Groups[4].Value
$splitAt   = $text.IndexOf('__________________') + 2
$textA     = $text.Substring(0, $splitAt)
$textB     = $text.Substring($splitAt)

$replacement = "$runOpen$runProps$tOpen$textA</w:t></w:r>" +
               '<w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/>' +
               "$runOpen$runProps$tOpen$textB</w:t></w:r>"

$frag = $frag.Substring(0, $runMatch.Index) + $replacement + $frag.Substring($runMatch.Index + $runMatch.Length)
$pLine.Range.InsertXML($frag)

# ---------------------------------------------------------------------
# 6. "/Ф.И.О./ ... /подпись/" line: only the run carries the highlight.
# ---------------------------------------------------------------------
$pRoles = Get-ParaByIndex ($idxLabel + 3)
$frag = Get-ParaXml $pRoles
$frag = Remove-Highlight $frag
$pRoles.Range.InsertXML($frag)

Write-Output "done"
